# Leave Card update (12/22/2023 10:59 AM)
# - Insert a new table row before the old "A53" (date 45292) row, shifting
#   all subsequent leave-card rows down by one.
# - Turn the newly inserted row into a "2024" year-marker row (matching the
#   existing 2021/2022/2023 marker rows).
# - Record a new SL entitlement pair in rows 50 & 51 (SL(1-0-0), 1.25 earned,
#   1 day used, dated 10/9/2023 and 11/9/2023).
# - Record the already-existing SP(2-0-0) row (row 49) EARNED amount (1.25).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")
$tbl = $ws.ListObjects.Item(1)

# --- Insert a new physical row at sheet row 53 (pushes 53..136 -> 54..137) ---
$ws.Rows.Item(53).Insert()

# The Excel Table itself doesn't auto-grow from a plain row Insert in this
# host, so resize it explicitly to cover the new row.
$tbl.Resize($ws.Range("A8:K137"))

# Resizing (while the old bottom "totals" row formula got pushed to 137)
# leaves the calculated-column formula on G137 in its unqualified "@" form,
# which evaluates incorrectly outside of a table context - restore it to
# match the rest of the calculated column.
$ws.Range("G137").Formula = '=IF(ISBLANK(Table1[[#This Row],[EARNED]]),"",Table1[[#This Row],[EARNED]])'

# The freshly inserted row 53 has generic/default formatting - copy the
# normal data-row look from row 54 (directly below, still carrying the old
# row-53 formatting) so borders/number-formats/fonts match the table body.
$ws.Range("A54:K54").Copy()
$ws.Range("A53:K53").PasteSpecial(-4122)   # xlPasteFormats

# --- Data entry: rows 49-51 (new SP / SL leave usage) ---
$ws.Range("C49").Value = 1.25

$ws.Range("B50").Value = "SL(1-0-0)"
$ws.Range("C50").Value = 1.25
$ws.Range("H50").Value = 1

$ws.Range("B51").Value = "SL(1-0-0)"
$ws.Range("C51").Value = 1.25
$ws.Range("H51").Value = 1

# K50 / K51 carry dates (like K49) - copy K49's date formatting onto them
# before writing the serial date values.
$ws.Range("K49").Copy()
$ws.Range("K50").PasteSpecial(-4122)
$ws.Range("K51").PasteSpecial(-4122)
$ws.Range("K50").Value = 45208   # 10/9/2023
$ws.Range("K51").Value = 45239   # 11/9/2023

# --- Turn row 53 into the "2024" year-marker row (like rows 10/21/35) ---
$ws.Range("A53").NumberFormat = "@"
$ws.Range("A53").Value = "2024"
$ws.Range("A35").Copy()            # A35 = existing "2023" marker cell
$ws.Range("A53").PasteSpecial(-4122)

$excel.CutCopyMode = $false

# --- View-state tweaks from the diff (selection only, cosmetic) ---
$ws.Range("I9").Select()
$ws.Range("F48").Select()
